$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text in the sheet (inline strings), even when
# they look numeric (e.g. "594.39" or "72.322.91"). Force text interpretation
# via a temporary Text number format, then restore the default "Normal" style so
# no stray cell formatting is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "72.322.91"
$ws.Range("E2").Value = "  +4.48%  "
Set-TextValue $ws.Range("D3") "3.611.56"
$ws.Range("E3").Value = "  +6.87%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue $ws.Range("D5") "594.39"
$ws.Range("E5").Value = "  +1.37%  "
Set-TextValue $ws.Range("D6") "184.04"
$ws.Range("E6").Value = "  +2.92%  "
Set-TextValue $ws.Range("D7") "3.603.75"
$ws.Range("E7").Value = "  +6.82%  "
Set-TextValue $ws.Range("D8") "0.609"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +7.18%  "
$ws.Range("E11").Value = "  +3.65%  "
Set-TextValue $ws.Range("D12") "50.19"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("E13").Value = "  +4.37%  "
Set-TextValue $ws.Range("D14") "696.22"
$ws.Range("E14").Value = "  +1.74%  "
Set-TextValue $ws.Range("D15") "4.190.57"
$ws.Range("E15").Value = "  +6.89%  "
$ws.Range("E16").Value = "  +4.00%  "
Set-TextValue $ws.Range("D17") "72.322.46"
$ws.Range("E17").Value = "  +4.43%  "
Set-TextValue $ws.Range("D18") "3.573.10"
$ws.Range("E18").Value = "  +5.75%  "
$ws.Range("E19").Value = "  +1.57%  "
Set-TextValue $ws.Range("D20") "18.52"
$ws.Range("E20").Value = "  +5.11%  "
Set-TextValue $ws.Range("D21") "11.74"
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("E22").Value = "  +3.58%  "
Set-TextValue $ws.Range("D23") "5.72"
$ws.Range("E23").Value = "  +5.61%  "
Set-TextValue $ws.Range("D24") "17.75"
$ws.Range("E24").Value = "  +3.43%  "
Set-TextValue $ws.Range("D25") "104.91"
$ws.Range("E25").Value = "  +1.63%  "
Set-TextValue $ws.Range("D26") "4.02"
$ws.Range("E26").Value = "  +2.57%  "
Set-TextValue $ws.Range("D27") "2.85"
$ws.Range("E27").Value = "  +4.81%  "
Set-TextValue $ws.Range("D28") "10.10"
$ws.Range("E28").Value = "  +5.00%  "
Set-TextValue $ws.Range("D29") "35.30"
$ws.Range("E29").Value = "  +4.21%  "
Set-TextValue $ws.Range("D30") "9.07"
Set-TextValue $ws.Range("D31") "7.54"
$ws.Range("E31").Value = "  +8.67%  "
$ws.Range("E32").Value = "  +16.81%  "
Set-TextValue $ws.Range("D33") "595.12"
$ws.Range("E33").Value = "  +6.17%  "
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  +1.31%  "
Set-TextValue $ws.Range("D36") "60.03"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("E37").Value = "  +0.06%  "
Set-TextValue $ws.Range("D38") "3.676.61"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  +5.04%  "
Set-TextValue $ws.Range("D40") "0.0₃0788"
$ws.Range("E40").Value = "  +13.34%  "
Set-TextValue $ws.Range("D41") "36.30"
$ws.Range("E41").Value = "  +1.62%  "
Set-TextValue $ws.Range("D42") "3.49"
$ws.Range("E42").Value = "  +7.35%  "
Set-TextValue $ws.Range("D43") "2.84"
$ws.Range("E43").Value = "  +6.02%  "
Set-TextValue $ws.Range("D44") "0.0442"
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  +4.13%  "
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("E48").Value = "  +4.55%  "
Set-TextValue $ws.Range("D49") "0.132"
$ws.Range("E49").Value = "  +2.04%  "
Set-TextValue $ws.Range("D50") "1.00"
$ws.Range("E50").Value = "  -0.09%  "
Set-TextValue $ws.Range("D51") "133.84"
$ws.Range("E51").Value = "  +0.46%  "
